$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5560.25
$ws.Range("I51").Value = 5333.3335
$ws.Range("J51").Value = 5696.4
$ws.Range("K51").Value = 5333.3335
$ws.Range("L51").Value = 5696.4
$ws.Range("M51").Value = -4849.3335
$ws.Range("N51").Value = -6664.4

$ws.Range("H86").Value = 90913224
$ws.Range("I86").Value = 10003
$ws.Range("J86").Value = 100003550
$ws.Range("K86").Value = 10003
$ws.Range("L86").Value = 100003550
$ws.Range("M86").Value = -8880
$ws.Range("N86").Value = -100005796

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()

$ws.Range("H89").Value = 90913224
$ws.Range("I89").Value = 10003
$ws.Range("J89").Value = 100003550
$ws.Range("K89").Value = 50015
$ws.Range("L89").Value = 500017750
$ws.Range("M89").Value = -44399
$ws.Range("N89").Value = -500028982

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()

$ws.Range("H132").Value = 1526.9667
$ws.Range("I132").Value = 1612.875
$ws.Range("J132").Value = 1183.3334
$ws.Range("K132").Value = 4838.625
$ws.Range("L132").Value = 3550.0002
$ws.Range("M132").Value = -2308.625
$ws.Range("N132").Value = -8610.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 1075
$ws.Range("I31").Value = 1075
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1075
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -781
$ws.Range("N31").ClearContents()

$ws.Range("H32").Value = 34862.523
$ws.Range("I32").Value = 32241.945
$ws.Range("K32").Value = 32241.945
$ws.Range("M32").Value = -31954.945

$ws.Range("H37").Value = 9024.799999999999
$ws.Range("I37").Value = 5016
$ws.Range("J37").Value = 15038
$ws.Range("K37").Value = 5016
$ws.Range("L37").Value = 15038
$ws.Range("M37").Value = -4743
$ws.Range("N37").Value = -15584

$ws.Range("H44").Value = 8888
$ws.Range("J44").Value = 8888
$ws.Range("L44").Value = 8888
$ws.Range("N44").Value = -9864

$ws.Range("H55").Value = 14852
$ws.Range("J55").Value = 15822.4
$ws.Range("L55").Value = 15822.4
$ws.Range("N55").Value = -16452.4

$ws.Range("H74").Value = 802.6957
$ws.Range("I74").Value = 809.6
$ws.Range("J74").Value = 789.75
$ws.Range("K74").Value = 809.6
$ws.Range("L74").Value = 789.75
$ws.Range("M74").Value = 64.39999999999998
$ws.Range("N74").Value = -2537.75

$ws.Range("H77").Value = 802.6957
$ws.Range("I77").Value = 809.6
$ws.Range("J77").Value = 789.75
$ws.Range("K77").Value = 4048
$ws.Range("L77").Value = 3948.75
$ws.Range("M77").Value = 320
$ws.Range("N77").Value = -12684.75

$ws.Range("H80").Value = 17532.5
$ws.Range("J80").Value = 20110
$ws.Range("L80").Value = 20110
$ws.Range("N80").Value = -22106

$ws.Range("H83").Value = 17532.5
$ws.Range("J83").Value = 20110
$ws.Range("L83").Value = 60330
$ws.Range("N83").Value = -70314

$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3200
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 3200
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -1578
$ws.Range("N102").Value = -5244

$ws.Range("H122").Value = 1570.6666
$ws.Range("I122").Value = 1570.6666
$ws.Range("K122").Value = 4711.9998
$ws.Range("M122").Value = -2261.9998

$ws.Range("H132").Value = 2933.2222
$ws.Range("I132").Value = 2198.8572
$ws.Range("J132").Value = 3400.5454
$ws.Range("K132").Value = 6596.571599999999
$ws.Range("L132").Value = 10201.6362
$ws.Range("M132").Value = -4066.571599999999
$ws.Range("N132").Value = -15261.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2153.3809
$ws.Range("I86").Value = 1930.8235
$ws.Range("K86").Value = 1930.8235
$ws.Range("M86").Value = -807.8235

$ws.Range("H89").Value = 2153.3809
$ws.Range("I89").Value = 1930.8235
$ws.Range("K89").Value = 9654.1175
$ws.Range("M89").Value = -4038.1175

$ws.Range("H105").Value = 2891.2144
$ws.Range("I105").Value = 2887.5
$ws.Range("J105").Value = 2897.9
$ws.Range("K105").Value = 2887.5
$ws.Range("L105").Value = 2897.9
$ws.Range("M105").Value = -1140.5
$ws.Range("N105").Value = -6391.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 293.5
$ws.Range("I22").Value = 267.75
$ws.Range("J22").Value = 345
$ws.Range("K22").Value = 267.75
$ws.Range("L22").Value = 345
$ws.Range("M22").Value = 82.25
$ws.Range("N22").Value = -1045

$ws.Range("H31").Value = 1566.2413
$ws.Range("I31").Value = 1248.48
$ws.Range("K31").Value = 1248.48
$ws.Range("M31").Value = -953.48

$ws.Range("H32").Value = 3000505
$ws.Range("I32").Value = 3000505
$ws.Range("K32").Value = 3000505
$ws.Range("M32").Value = -3000189

$ws.Range("H34").Value = 1566.2413
$ws.Range("I34").Value = 1248.48
$ws.Range("K34").Value = 1248.48
$ws.Range("M34").Value = -1046.48

$ws.Range("H35").Value = 700
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H132").Value = 2410.25
$ws.Range("I132").Value = 1559.6666
$ws.Range("K132").Value = 4678.9998
$ws.Range("M132").Value = -2148.9998

$ws.Range("H140").Value = 51705.26
$ws.Range("J140").Value = 51705.26
$ws.Range("L140").Value = 51705.26
$ws.Range("N140").Value = -62065.26

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 30003.5
$ws.Range("I29").Value = 40007
$ws.Range("J29").Value = 20000
$ws.Range("K29").Value = 40007
$ws.Range("L29").Value = 20000
$ws.Range("M29").Value = -39717
$ws.Range("N29").Value = -20580

$ws.Range("H80").Value = 3666.6667
$ws.Range("I80").Value = 3666.6667
$ws.Range("K80").Value = 3666.6667
$ws.Range("M80").Value = -2668.6667

$ws.Range("H83").Value = 3666.6667
$ws.Range("I83").Value = 3666.6667
$ws.Range("K83").Value = 18333.3335
$ws.Range("M83").Value = -13341.3335

$ws.Range("H132").Value = 2689.5
$ws.Range("I132").Value = 2083.2
$ws.Range("J132").Value = 3555.6428
$ws.Range("K132").Value = 6249.599999999999
$ws.Range("L132").Value = 10666.9284
$ws.Range("M132").Value = -3719.599999999999
$ws.Range("N132").Value = -15726.9284

$ws.Range("H138").Value = 19800
$ws.Range("J138").Value = 19800
$ws.Range("L138").Value = 19800
$ws.Range("N138").Value = -30080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 500
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1090

$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 500
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -714

$ws.Range("H32").Value = 50000000
$ws.Range("I32").Value = 50000000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 50000000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -49999683
$ws.Range("N32").ClearContents()

$ws.Range("H46").Value = 467
$ws.Range("I46").Value = 467
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 467
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -279
$ws.Range("N46").ClearContents()

$ws.Range("H132").Value = 4566.8
$ws.Range("I132").Value = 4286.857
$ws.Range("J132").Value = 4811.75
$ws.Range("K132").Value = 12860.571
$ws.Range("L132").Value = 14435.25
$ws.Range("M132").Value = -10330.571
$ws.Range("N132").Value = -19495.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 5005002.5
$ws.Range("I13").Value = 5005002.5
$ws.Range("K13").Value = 5005002.5
$ws.Range("M13").Value = -5004862.5

$ws.Range("H132").Value = 1241.2258
$ws.Range("I132").Value = 658.75
$ws.Range("J132").Value = 1862.5333
$ws.Range("K132").Value = 1976.25
$ws.Range("L132").Value = 5587.5999
$ws.Range("M132").Value = 553.75
$ws.Range("N132").Value = -10647.5999

$ws.Range("H136").Value = 2494.652
$ws.Range("I136").Value = 3208.2964
$ws.Range("J136").Value = 1480.5264
$ws.Range("K136").Value = 9624.889200000001
$ws.Range("L136").Value = 4441.5792
$ws.Range("M136").Value = -7074.889200000001
$ws.Range("N136").Value = -9541.5792
